$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (plain text, no numeric risk) ---
$ws.Range("B35").Value = "VeChain"
$ws.Range("B36").Value = "Hedera"
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("B42").Value = "Frax"
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("B44").Value = "FraxShare"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("B47").Value = "Decentraland"
$ws.Range("B48").Value = "Quant"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("B50").Value = "EOS"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

# --- Price column (D): force text to preserve literal formatting ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.938.49"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.718.39"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.11"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4867"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3499"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41.99"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07249"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.046"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.88"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.718.67"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.854"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.54"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001038"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06364"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.0000"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.49"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.645"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "26.995.98"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.78"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.081"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.56"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.99"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.911.12"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.073"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.97"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.027"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09291"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.581"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.348"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02180"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05879"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.452"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1993"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6013"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.730"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9992"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.094"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.509"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.70"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.572"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5632"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.36"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.832"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.109"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06652"

# --- Volume(1h) column (E): padded percentage strings stay text naturally ---
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E5").Value = "  -6.05%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  +7.13%  "
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -4.58%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("E16").Value = "  -4.63%  "
$ws.Range("E17").Value = "  -6.35%  "
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("E26").Value = "  -5.39%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("E29").Value = "  -3.76%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  -4.07%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("E35").Value = "  -4.25%  "
$ws.Range("E36").Value = "  -4.03%  "
$ws.Range("E37").Value = "  +5.14%  "
$ws.Range("E38").Value = "  -7.19%  "
$ws.Range("E39").Value = "  -4.41%  "
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("E41").Value = "  -4.21%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  -7.43%  "
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("E46").Value = "  -4.36%  "
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("E49").Value = "  -5.33%  "
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("E51").Value = "  -2.33%  "
